$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the shared text used by the "MODEL_CONDITION" header to "MODELCONDITION"
# (same cell, same position -- text-only rename).
$ws.Cells.Replace("MODEL_CONDITION", "MODELCONDITION") | Out-Null

# Drop the old leading index column -- the remaining columns shift one
# place to the left (A1:F3 -> A1:E3).
$ws.Columns("A").Delete()
